# Financials update: insert a new "most recent period" column before column D
# on the CWK sheet, shifting the existing D:K data right to E:L, and populate
# the new column D with the latest figures (per the commit's data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before D; existing D:K shift to E:L automatically
#    (Excel also extends the trailing "NA" markers into the new last column).
$ws.Columns("D").Insert()

# 2) The new column D currently inherits formatting from column C (to its
#    left). Re-apply the correct number formats/styles by copying them from
#    column E (the former column D, now shifted one to the right).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# 3) Populate new column D with the latest-period values for rows that carry
#    real figures. (Rows not listed here are either blank spacer rows, the
#    all-"NA" row, or all-zero rows, which already read correctly after the
#    insert above and need no further edit.)
$newValues = @{
    7 = 43465
    8 = 8219900
    9 = 6642400
    10 = 1577500
    14 = 3800
    15 = 290000
    17 = 8207300
    18 = 12600
    20 = 15000
    21 = 317600
    22 = 238400
    23 = -210800
    24 = 4200
    26 = -215000
    27 = -215000
    29 = 29200
    32 = -15000
    33 = -185800
    35 = -185800
    38 = 43465
    41 = 895300
    43 = 1665200
    45 = 182800
    46 = 2743300
    47 = 9200
    48 = 293600
    49 = 2926900
    52 = 573000
    54 = 6546000
    57 = 980900
    58 = 39900
    59 = 1017900
    60 = 2038700
    61 = 2644200
    62 = 503000
    66 = 5185900
    72 = -1298400
    76 = 1360100
    80 = 43465
    81 = -185800
    83 = 290000
    89 = -2200
    91 = -84200
    94 = -218000
    100 = 725900
    101 = -8200
    102 = 497500
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 4).Value2 = $newValues[$row]
}

# 4) Rows whose whole data band reads "0" across D:K keep that same 0 in the
#    new column D (content is unchanged in appearance, just shifted right).
$zeroRows = @(13,25,28,30,31,34,42,44,50,51,53,63,64,65,68,69,70,71,73,74,75,77,84,85,86,87,88,92,93,96,97,98,99)
foreach ($row in $zeroRows) {
    $ws.Cells.Item($row, 4).Value2 = 0
}

# 5) Row 12 reads "NA" across the whole D:K band; keep that in new column D
#    too (copy the text straight from the neighboring cell to reuse the
#    existing shared string instead of minting a new one).
$ws.Range("D12").Value2 = $ws.Range("E12").Value2
